$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the shared-string header in T1 from "%dR_c" to "%dR_u"
$ws.Range("T1").Value = "%dR_u"

# 2. Change the R-model uncertainty values in column T (rows 2-16) from 2 (2%, Simona's
#    paper) to 0.7 (0.7%, Christy)
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 20).Value = 0.7
}

# 3. Widen columns A:T (and beyond, matching the original column definitions) slightly
$ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(1,1025)).EntireColumn.ColumnWidth = 13.5

# 4. Move the active selection to F9
$ws.Range("F9").Select()
